$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.461.78'
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.898.79'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.01'
$ws.Range("E5").Value = '  +0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4911'
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.899.78'
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.95'
$ws.Range("E11").Value = '  +1.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07331'
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.171'
$ws.Range("E13").Value = '  +3.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.64'
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6654'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.440.27'
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.46'
$ws.Range("E17").Value = '  +3.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007838'
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9996'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.144.78'
$ws.Range("E20").Value = '  +0.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.357'
$ws.Range("E21").Value = '  +13.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '192.82'
$ws.Range("E23").Value = '  +1.45%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.110'
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.493'
$ws.Range("E25").Value = '  +2.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.72'
$ws.Range("E26").Value = '  +2.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.26'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.946'
$ws.Range("E28").Value = '  +6.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.470'
$ws.Range("E29").Value = '  +4.65%  '
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09195'
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.047'
$ws.Range("E32").Value = '  +3.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05195'
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7419'
$ws.Range("E34").Value = '  +2.35%  '
$ws.Range("E35").Value = '  +1.90%  '
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.676'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9239'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.042'
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4385'
$ws.Range("E41").Value = '  +0.23%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.912'
$ws.Range("E42").Value = '  +3.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.32'
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9946'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.23'
$ws.Range("E45").Value = '  +21.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1371'
$ws.Range("E46").Value = '  +3.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.615'
$ws.Range("E47").Value = '  +4.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.031'
$ws.Range("E48").Value = '  +4.62%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.05'
$ws.Range("E49").Value = '  +5.45%  '
$ws.Range("E50").Value = '  -0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3926'
$ws.Range("E51").Value = '  -2.86%  '

Write-Output "Updated cryptos list"
